# The workbook's loan product code was missing a hyphen after "343".
# Fix the product-code value stored in cell B1 on both sheets, then
# leave the workbook with "ProductLoanOutput" selected as the active
# tab/sheet (with B1 selected there), matching how it was re-saved.

$wb = $excel.ActiveWorkbook

$wsInput  = $wb.Worksheets.Item(1)   # "ProductLoanInput"
$wsOutput = $wb.Worksheets.Item(2)   # "ProductLoanOutput"

$newCode = "343-MS-EPP-DB-SAR-REC-NON-RNI-CTRFD-SAR-MD-TR-1-ONTIME"

$wsInput.Range("B1").Value = $newCode
$wsOutput.Range("B1").Value = $newCode

# Select B1 on the input sheet (was previously A6:B6).
$wsInput.Activate()
$wsInput.Range("B1").Select()

# Finish with the output sheet active/selected at B1 - this becomes
# the workbook's active tab.
$wsOutput.Activate()
$wsOutput.Range("B1").Select()
